$d = $word.ActiveDocument

function ReplaceAll($find, $replace) {
    $r = $d.Content
    $r.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# Row 5 ("point revisi: Small figure." / "hasil revisi: Made the figure bigger.")
# collapse the split runs into single runs with the same text.
ReplaceAll "point revisi: Small figure." "point revisi: Small figure."
ReplaceAll "hasil revisi: Made the figure bigger." "hasil revisi: Made the figure bigger."

# Row 6 ("point revisi: Add general flowchart on the system." / "hasil revisi: Added general flowchart before sub-chapters.")
ReplaceAll "point revisi: Add general flowchart on the system." "point revisi: Add general flowchart on the system."
ReplaceAll "hasil revisi: Added general flowchart before sub-chapters." "hasil revisi: Added general flowchart before sub-chapters."

# Row 7 ("point revisi: Summary table on datasets." / "hasil revisi: Created summary table on dataset.")
ReplaceAll "point revisi: Summary table on datasets." "point revisi: Summary table on datasets."
ReplaceAll "hasil revisi: Created summary table on dataset." "hasil revisi: Created summary table on dataset."

# Row 8 ("point revisi: Show each step in pre-processing step for one data." / "hasil revisi: Added explanation on pre-processing step for one data.")
ReplaceAll "point revisi: Show each step in pre-processing step for one data." "point revisi: Show each step in pre-processing step for one data."
ReplaceAll "hasil revisi: Added explanation on pre-processing step for one data." "hasil revisi: Added explanation on pre-processing step for one data."

# Row 10 ("point revisi: Change conclusion points to numbers." / "hasil revisi: Changed conclusion points to numbers.")
ReplaceAll "point revisi: Change conclusion points to numbers." "point revisi: Change conclusion points to numbers."
ReplaceAll "hasil revisi: Changed conclusion points to numbers." "hasil revisi: Changed conclusion points to numbers."

# Closing date line: "Surabaya, 23 Juni 2023" -> "Surabaya, 4 Juli 2023"
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Surabaya*23 Juni*") {
        $rng = $p.Range.Duplicate
        $rng.Find.Execute("23 Juni", $true, $false, $false, $false, $false, $true, 1, $false, "4 Juli", 2) | Out-Null
    }
}
